$wb = $excel.ActiveWorkbook

$oldGuid = "285d2b4f-c17a-4342-84ac-1c63d0f35aa6"
$newGuid = "f0a3bbd5-d330-4d1d-9c39-132ba6ce4c02"

$oldZhXlf = "285d2b4f-c17a-4342-84ac-1c63d0f35aa6.625348eb0cedbb4f26c27554be30eed7f36c7f65.zh-cn.xlf"
$newZhXlf = "f0a3bbd5-d330-4d1d-9c39-132ba6ce4c02.8c177846d09af2314b37419600aa30511e2232ad.zh-cn.xlf"

$oldDeXlf = "285d2b4f-c17a-4342-84ac-1c63d0f35aa6.625348eb0cedbb4f26c27554be30eed7f36c7f65.de-de.xlf"
$newDeXlf = "f0a3bbd5-d330-4d1d-9c39-132ba6ce4c02.8c177846d09af2314b37419600aa30511e2232ad.de-de.xlf"

$zeroDate = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = $newGuid + ".md"

$newPathAndName = "e2e\" + $newGuid + ".md"
$ov.Range("B2").Value = $newPathAndName
foreach ($h in $ov.Hyperlinks) {
    if ($h.Range.Address() -eq '$B$2') {
        $h.TextToDisplay = $newPathAndName
    }
}

$ov.Range("G2").Value = "2016-09-04 17:05:37"

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = $newGuid + ".md"
foreach ($h in $zh.Hyperlinks) {
    if ($h.Range.Address() -eq '$A$2') {
        $h.TextToDisplay = $newGuid + ".md"
    }
}

$zh.Range("G2").Value = $newZhXlf
$zh.Range("H2").Value = "2016-09-04 17:05:32"

# Latest Target File (I2) loses its value, style and hyperlink entirely.
$zhKeepAddr = $null
foreach ($h in $zh.Hyperlinks) {
    if ($h.Range.Address() -eq '$A$2') {
        $zhKeepAddr = $h.Address
    }
}
$zh.Range("I2").Value = ""
$zh.Range("I2").ClearFormats()
$zh.Range("I2").Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), $zhKeepAddr, "", "", $newGuid + ".md") | Out-Null

$zh.Range("J2").Value = ""
$zh.Range("K2").Value = $zeroDate

$zh.Columns.Item(9).ColumnWidth = 18.6506053379604
$zh.Columns.Item(10).ColumnWidth = 21.7054770333426

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = $newGuid + ".md"

$de.Range("G2").Value = $newDeXlf
$de.Range("H2").Value = "2016-09-04 17:05:37"

$deKeepAddr = $null
foreach ($h in $de.Hyperlinks) {
    if ($h.Range.Address() -eq '$A$2') {
        $deKeepAddr = $h.Address
    }
}
$de.Range("I2").Value = ""
$de.Range("I2").ClearFormats()
$de.Range("I2").Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), $deKeepAddr, "", "", $newGuid + ".md") | Out-Null

$de.Range("J2").Value = ""
$de.Range("K2").Value = $zeroDate

$de.Columns.Item(9).ColumnWidth = 18.6506053379604
$de.Columns.Item(10).ColumnWidth = 21.7054770333426
